# Update "想去人数" (F column) figures on the 展览 (sheet1) and 全部类型 (sheet4)
# worksheets, matching the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitChanges = @{
    4  = 8073
    5  = 102
    7  = 1284
    8  = 33
    10 = 486
    11 = 169
    14 = 72
    15 = 82
    17 = 5982
    18 = 194
    19 = 286
    20 = 2078
    21 = 58
    22 = 93
    23 = 242
    24 = 424
}
foreach ($row in $exhibitChanges.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitChanges[$row]
}

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$allChanges = @{
    4  = 8073
    5  = 102
    7  = 1284
    8  = 33
    11 = 486
    12 = 169
    15 = 72
    16 = 82
    19 = 5982
    21 = 194
    22 = 286
    23 = 2079
    24 = 58
    25 = 93
    26 = 242
    27 = 424
}
foreach ($row in $allChanges.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allChanges[$row]
}
